$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 892; this shifts the existing rows 892..933
# down to 893..934 and extends the used range to row 934.
$ws.Rows.Item(892).Insert()

# Column A holds a "yyyy/mm/dd" string that must stay plain text (the
# sheet stores dates as literal strings, not Excel date serials).
# Assigning the literal text straight into .Value triggers Excel's
# autodetection of date-looking strings and turns the cell into a real
# date serial with a date number format, which the source file does not
# have. Writing it first as a formula that evaluates to a text string,
# then collapsing that formula down to a static value via Copy /
# PasteSpecial(xlPasteValues), keeps it a plain text cell with no
# number-format override.
$ws.Cells.Item(892, 1).Formula = '="2026/02/27"'
$ws.Cells.Item(892, 1).Copy()
$ws.Cells.Item(892, 1).PasteSpecial(-4163)

$ws.Cells.Item(892, 2).Value = "金"
$ws.Cells.Item(892, 3).Value = 1
$ws.Cells.Item(892, 4).Value = 201
